$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B2").Value = 25
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 21
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 8
$ws.Range("B7").Value = 9
$ws.Range("B8").Value = 12
$ws.Range("B9").Value = 20
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 23
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 7
$ws.Range("B15").Value = 11

$ws.Range("C2").Value = "P.J. Washington"
$ws.Range("C3").Value = "Terry Rozier"
$ws.Range("C4").Value = "JT Thor"
$ws.Range("C5").Value = "Nick Richards"
$ws.Range("C6").Value = "Dennis Smith Jr."
$ws.Range("C7").Value = "Théo Maledon (TW)"
$ws.Range("C8").Value = "Kelly Oubre Jr."
$ws.Range("C9").Value = "Gordon Hayward"
$ws.Range("C10").Value = "LaMelo Ball"
$ws.Range("C11").Value = "Mark Williams"
$ws.Range("C12").Value = "Kai Jones"
$ws.Range("C13").Value = "James Bouknight"
$ws.Range("C14").Value = "Bryce McGowens (TW)"
$ws.Range("C15").Value = "Cody Martin"
$ws.Range("C16").Value = "Reggie Jackson"
$ws.Range("C17").Value = "Svi Mykhailiuk"

$ws.Range("D2").Value = "PF"
$ws.Range("D3").Value = "SG"
$ws.Range("D4").Value = "PF"
$ws.Range("D5").Value = "C"
$ws.Range("D6").Value = "PG"
$ws.Range("D7").Value = "PG"
$ws.Range("D8").Value = "SF"
$ws.Range("D9").Value = "SF"
$ws.Range("D10").Value = "PG"
$ws.Range("D11").Value = "C"
$ws.Range("D12").Value = "C"
$ws.Range("D13").Value = "SG"
$ws.Range("D14").Value = "SG"
$ws.Range("D15").Value = "SF"
$ws.Range("D16").Value = "PG"
$ws.Range("D17").Value = "SF"

$ws.Range("E2").Value = "6-7"
$ws.Range("E3").Value = "6-1"
$ws.Range("E4").Value = "6-10"
$ws.Range("E5").Value = "7-0"
$ws.Range("E6").Value = "6-2"
$ws.Range("E7").Value = "6-4"
$ws.Range("E8").Value = "6-7"
$ws.Range("E9").Value = "6-7"
$ws.Range("E10").Value = "6-7"
$ws.Range("E11").Value = "7-1"
$ws.Range("E12").Value = "6-11"
$ws.Range("E13").Value = "6-5"
$ws.Range("E14").Value = "6-7"
$ws.Range("E15").Value = "6-5"
$ws.Range("E16").Value = "6-2"
$ws.Range("E17").Value = "6-7"

$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 190
$ws.Range("F4").Value = 205
$ws.Range("F5").Value = 245
$ws.Range("F6").Value = 205
$ws.Range("F7").Value = 175
$ws.Range("F8").Value = 203
$ws.Range("F9").Value = 225
$ws.Range("F10").Value = 180
$ws.Range("F11").Value = 241
$ws.Range("F12").Value = 218
$ws.Range("F13").Value = 190
$ws.Range("F14").Value = 179
$ws.Range("F15").Value = 205
$ws.Range("F16").Value = 208
$ws.Range("F17").Value = 205

$ws.Range("G2").Value = "August 23, 1998"
$ws.Range("G3").Value = "March 17, 1994"
$ws.Range("G4").Value = "August 26, 2002"
$ws.Range("G5").Value = "November 29, 1997"
$ws.Range("G6").Value = "November 25, 1997"
$ws.Range("G7").Value = "June 12, 2001"
$ws.Range("G8").Value = "December 9, 1995"
$ws.Range("G9").Value = "March 23, 1990"
$ws.Range("G10").Value = "August 22, 2001"
$ws.Range("G11").Value = "December 16, 2001"
$ws.Range("G12").Value = "January 19, 2001"
$ws.Range("G13").Value = "September 18, 2000"
$ws.Range("G14").Value = "November 8, 2002"
$ws.Range("G15").Value = "September 28, 1995"
$ws.Range("G16").Value = "April 16, 1990"
$ws.Range("G17").Value = "June 10, 1997"

$ws.Range("H2").Value = "us"
$ws.Range("H3").Value = "us"
$ws.Range("H4").Value = "us"
$ws.Range("H5").Value = "jm"
$ws.Range("H6").Value = "us"
$ws.Range("H7").Value = "fr"
$ws.Range("H8").Value = "us"
$ws.Range("H9").Value = "us"
$ws.Range("H10").Value = "us"
$ws.Range("H11").Value = "us"
$ws.Range("H12").Value = "bs"
$ws.Range("H13").Value = "us"
$ws.Range("H14").Value = "us"
$ws.Range("H15").Value = "us"
$ws.Range("H16").Value = "it"
$ws.Range("H17").Value = "ua"

$ws.Range("I2").Value = "3"
$ws.Range("I3").Value = "7"
$ws.Range("I4").Value = "1"
$ws.Range("I5").Value = "2"
$ws.Range("I6").Value = "5"
$ws.Range("I7").Value = "2"
$ws.Range("I8").Value = "7"
$ws.Range("I9").Value = "12"
$ws.Range("I10").Value = "2"
$ws.Range("I11").Value = "R"
$ws.Range("I12").Value = "1"
$ws.Range("I13").Value = "1"
$ws.Range("I14").Value = "R"
$ws.Range("I15").Value = "3"
$ws.Range("I16").Value = "11"
$ws.Range("I17").Value = "4"

$ws.Range("J10").ClearContents()
$ws.Range("J2").Value = "Kentucky"
$ws.Range("J3").Value = "Louisville"
$ws.Range("J4").Value = "Auburn"
$ws.Range("J5").Value = "Kentucky"
$ws.Range("J6").Value = "NC State"
$ws.Range("J8").Value = "Kansas"
$ws.Range("J9").Value = "Butler"
$ws.Range("J11").Value = "Duke"
$ws.Range("J12").Value = "Texas"
$ws.Range("J13").Value = "UConn"
$ws.Range("J14").Value = "Nebraska"
$ws.Range("J15").Value = "NC State, Nevada"
$ws.Range("J16").Value = "Boston College"
$ws.Range("J17").Value = "Kansas"

$ws.Range("K2").Value = "https://www.basketball-reference.com/players/w/washipj01.html"
$ws.Range("K3").Value = "https://www.basketball-reference.com/players/r/roziete01.html"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/t/thorjt01.html"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/r/richani01.html"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/s/smithde03.html"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/m/maledth01.html"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/o/oubreke01.html"
$ws.Range("K9").Value = "https://www.basketball-reference.com/players/h/haywago01.html"
$ws.Range("K10").Value = "https://www.basketball-reference.com/players/b/ballla01.html"
$ws.Range("K11").Value = "https://www.basketball-reference.com/players/w/willima07.html"
$ws.Range("K12").Value = "https://www.basketball-reference.com/players/j/joneska01.html"
$ws.Range("K13").Value = "https://www.basketball-reference.com/players/b/bouknja01.html"
$ws.Range("K14").Value = "https://www.basketball-reference.com/players/m/mcgowbr01.html"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/m/martico01.html"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/j/jacksre01.html"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/m/mykhasv01.html"
